$d = $word.ActiveDocument

$replacements = @(
    @("2025-10-25 Saturday", "2025-10-26 Sunday"),
    @("72×21=1512", "52×27=1404"),
    @("62×22=1364", "17×72=1224"),
    @("74×48=3552", "93×68=6324"),
    @("33×47=1551", "91×50=4550"),
    @("65×53=3445", "26×66=1716"),
    @("79×54=4266", "56×95=5320"),
    @("21×90=1890", "69×26=1794"),
    @("22×87=1914", "25×62=1550"),
    @("64×74=4736", "94×92=8648"),
    @("58×98=5684", "17×17=289"),
    @("61×39=2379", "73×91=6643"),
    @("94×58=5452", "21×16=336"),
    @("87×28=2436", "13×21=273"),
    @("36×46=1656", "18×70=1260"),
    @("65×30=1950", "54×18=972"),
    @("88×60=5280", "59×56=3304"),
    @("65×57=3705", "90×93=8370"),
    @("82×81=6642", "12×82=984"),
    @("87×47=4089", "69×75=5175"),
    @("24×61=1464", "65×26=1690"),
    @("82×56=4592", "76×38=2888"),
    @("85×15=1275", "52×93=4836"),
    @("33×25=825", "76×82=6232"),
    @("12×36=432", "85×87=7395"),
    @("89×90=8010", "90×62=5580")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
